# Scen_NCAP_NUC.xlsx edit script
# - Scrolls/selects the sheet view to E84 (topLeftCell A42)
# - Changes Pset_PN from NCAP_BND to CAP_BND on several rows (C column)
# - Updates E24:E29 values (and clears their number-format style)
# - Re-labels several "UP" (LimType) rows as "\I: " header rows (B column),
#   applying the same look (style 10) used by the other section headers

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Sheet view: scroll to row 42 and select E84
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 42
$win.ScrollColumn = 1
$ws.Range("E84").Select()

# ---------------------------------------------------------------------
# 2. C15:C20, C24:C29, C42:C47, C60:C65, C69:C74 : NCAP_BND -> CAP_BND
# ---------------------------------------------------------------------
$capBndCells = @(
    "C15","C16","C17","C18","C19","C20",
    "C24","C25","C26","C27","C28","C29",
    "C42","C43","C44","C45","C46","C47",
    "C60","C61","C62","C63","C64","C65",
    "C69","C70","C71","C72","C73","C74"
)
foreach ($addr in $capBndCells) {
    $ws.Range($addr).Value2 = "CAP_BND"
}

# ---------------------------------------------------------------------
# 3. E24:E29 new values, with the shaded "s=20" format cleared back to
#    the plain/default look (copy format from an unstyled cell, D78)
# ---------------------------------------------------------------------
$plainFormatSource = $ws.Range("D78")
$plainFormatSource.Copy()
$ws.Range("E24:E29").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("E24").Value2 = 0.61
$ws.Range("E25").Value2 = 1.55
$ws.Range("E26").Value2 = 6.85
$ws.Range("E27").Value2 = 9.23
$ws.Range("E28").Value2 = 17.77
$ws.Range("E29").Value2 = 20.86

# ---------------------------------------------------------------------
# 4. B78:B84 and B86:B92 : "UP" -> "\I: ", styled like other section
#    header rows (copy formatting from B77, which already carries it)
# ---------------------------------------------------------------------
$headerFormatSource = $ws.Range("B77")
$headerFormatSource.Copy()

$headerRows = @("B78","B79","B80","B81","B82","B83","B84","B86","B87","B88","B89","B90","B91","B92")
foreach ($addr in $headerRows) {
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

foreach ($addr in $headerRows) {
    $ws.Range($addr).Value2 = "\I: "
}
